$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in column L
$ws.Range("L1").Value = "break_on_off"

# Values for L2:L73 (break_on_off flags) - default 0, with 1s at specific rows
$values = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $values[$i]
}

# Update selection to mirror the authored view state
$ws.Range("L1:L73").Select()
